$d = $word.ActiveDocument

# Locate the empty "ListParagraph" paragraph that immediately follows
# "Watching a video from different IP addresses is important." --
# that's where the two new notes get appended.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd() -eq "Watching a video from different IP addresses is important.") {
        $targetIndex = $i + 1
        break
    }
}

$target = $d.Paragraphs.Item($targetIndex)
$r = $target.Range
$r.InsertParagraphAfter()

$notePara = $d.Paragraphs.Item($targetIndex + 1)
$notePara.Range.Text = "Note:   this approach seems to work"
$notePara.Style = "List Paragraph"

$notePara.Range.InsertParagraphAfter()
$secondPara = $d.Paragraphs.Item($targetIndex + 2)
$secondPara.Range.Text = "Parallely opening the URLs and clicking on the video"
$secondPara.Style = "List Paragraph"
